$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Excel alignment constants
$xlTop = -4160
$xlLeft = -4131

# Green fill used for the "Easy" difficulty tag (RGB 0,176,80 -> 00B050)
$easyFillColor = 5287936

# ---- Row 13 (new row): "Left View of Binary Tree" ----
$ws.Range("B13").Value = "GFG"
$ws.Range("B13").VerticalAlignment = $xlTop
$ws.Range("B13").HorizontalAlignment = $xlLeft
$ws.Range("B13").WrapText = $true

$ws.Range("C13").Value = "Left View of Binary Tree"

$ws.Range("D13").Value = "Java/Python"
$ws.Range("D13").VerticalAlignment = $xlTop

$ws.Range("E13").Value = "Easy"
$ws.Range("E13").Interior.Color = $easyFillColor

# ---- Row 14 (previously a mostly-empty row): "Right View of Binary Tree" ----
$ws.Range("B14").Value = 199

$ws.Range("C14").Value = "Right View of Binary Tree"

$ws.Range("D14").Value = "Java/Python"
$ws.Range("D14").VerticalAlignment = $xlTop

$ws.Range("E14").Value = "Easy"
$ws.Range("E14").Interior.Color = $easyFillColor

# ---- Move the active selection to E14 ----
$ws.Range("E14").Select()
